$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20-43 down to 21-44.
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new weekly data point.
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44629
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112031
$ws.Cells.Item(20, 7).Value = "Poroto verde"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 30000
$ws.Cells.Item(20, 12).Value = 30000
$ws.Cells.Item(20, 13).Value = 30000
$ws.Cells.Item(20, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 1200
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
